$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new date columns E1:P1
$headers = @{
    "E1" = "14/9/2022"
    "F1" = "15/9/2022"
    "G1" = "19/9/2022"
    "H1" = "20/9/2022"
    "I1" = "21/9/2022"
    "J1" = "22/9/2022"
    "K1" = "23/9/2022"
    "L1" = "24/9/2022"
    "M1" = "26/9/2022"
    "N1" = "27/9/2022"
    "O1" = "28/9/2022"
    "P1" = "29/9/2022"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

$ws.Range("D1").Copy()
$ws.Range("E1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-7 for columns E:P
$data = @{
    2 = @("não","sim","não","não","não","não","não","não","-","-","-","-")
    3 = @("não","sim","não","não","sim","não","não","não","não","sim","não","sim")
    4 = @("não","sim","não","não","não","não","não","não","-","-","-","-")
    5 = @("não","sim","não","não","sim","não","não","não","não","sim","não","sim")
    6 = @("não","sim","não","não","sim","não","não","não","não","sim","não","sim")
    7 = @("não","sim","não","não","não","não","não","não","-","-","-","-")
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = $values[$i]
    }
}
